# Auto-applied cell value updates per commit diff (Hyperion_Profits.xlsx)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 551.3333  # H28: 641.2 -> 551.3333
$ws.Cells.Item(28, 9).Value = 104  # I28: 104.5 -> 104
$ws.Cells.Item(28, 10).Value = 998.6667  # J28: 999 -> 998.6667
$ws.Cells.Item(28, 11).Value = 104  # K28: 104.5 -> 104
$ws.Cells.Item(28, 12).Value = 998.6667  # L28: 999 -> 998.6667
$ws.Cells.Item(28, 13).Value = 381  # M28: 380.5 -> 381
$ws.Cells.Item(28, 14).Value = -1968.6667  # N28: -1969 -> -1968.6667

$ws.Cells.Item(80, 8).Value = 3592.5  # H80: 3898.9092 -> 3592.5
$ws.Cells.Item(80, 9).Value = 405.5  # I80: 466.66666 -> 405.5
$ws.Cells.Item(80, 11).Value = 1216.5  # K80: 1399.99998 -> 1216.5
$ws.Cells.Item(80, 13).Value = -218.5  # M80: -401.9999800000001 -> -218.5

$ws.Cells.Item(83, 8).Value = 3592.5  # H83: 3898.9092 -> 3592.5
$ws.Cells.Item(83, 9).Value = 405.5  # I83: 466.66666 -> 405.5
$ws.Cells.Item(83, 11).Value = 3649.5  # K83: 4199.99994 -> 3649.5
$ws.Cells.Item(83, 13).Value = 1342.5  # M83: 792.0000600000003 -> 1342.5

$ws.Cells.Item(86, 8).Value = 4000  # H86: 4249.5 -> 4000
$ws.Cells.Item(86, 9).Value = 4000  # I86: 4249.5 -> 4000
$ws.Cells.Item(86, 11).Value = 4000  # K86: 4249.5 -> 4000
$ws.Cells.Item(86, 13).Value = -2877  # M86: -3126.5 -> -2877

$ws.Cells.Item(89, 8).Value = 4000  # H89: 4249.5 -> 4000
$ws.Cells.Item(89, 9).Value = 4000  # I89: 4249.5 -> 4000
$ws.Cells.Item(89, 11).Value = 20000  # K89: 21247.5 -> 20000
$ws.Cells.Item(89, 13).Value = -14384  # M89: -15631.5 -> -14384

$ws.Cells.Item(100, 8).Value = 5298  # H100: 7329.3335 -> 5298
$ws.Cells.Item(100, 9).Value = 5298  # I100: 7329.3335 -> 5298
$ws.Cells.Item(100, 11).Value = 5298  # K100: 7329.3335 -> 5298
$ws.Cells.Item(100, 13).Value = -4757  # M100: -6788.3335 -> -4757

$ws.Cells.Item(111, 8).Value = 15873734  # H111: 18519274 -> 15873734
$ws.Cells.Item(111, 9).Value = 22222828  # I111: 27778410 -> 22222828
$ws.Cells.Item(111, 11).Value = 66668484  # K111: 83335230 -> 66668484
$ws.Cells.Item(111, 13).Value = -66665417  # M111: -83332163 -> -66665417

$ws.Cells.Item(125, 8).Value = 6062926.5  # H125: 10103682 -> 6062926.5
$ws.Cells.Item(125, 9).Value = 1789.8572  # I125: 1949.5 -> 1789.8572
$ws.Cells.Item(125, 10).Value = 6946842.5  # J125: 12348511 -> 6946842.5
$ws.Cells.Item(125, 11).Value = 16108.7148  # K125: 17545.5 -> 16108.7148
$ws.Cells.Item(125, 12).Value = 62521582.5  # L125: 111136599 -> 62521582.5
$ws.Cells.Item(125, 13).Value = -13648.7148  # M125: -15085.5 -> -13648.7148
$ws.Cells.Item(125, 14).Value = -62526502.5  # N125: -111141519 -> -62526502.5

$ws.Cells.Item(137, 8).Value = 102276.39  # H137: 83780.77 -> 102276.39
$ws.Cells.Item(137, 9).Value = 359555  # I137: 224847.12 -> 359555
$ws.Cells.Item(137, 10).Value = 3323.077  # J137: 3171.4285 -> 3323.077
$ws.Cells.Item(137, 11).Value = 1078665  # K137: 674541.36 -> 1078665
$ws.Cells.Item(137, 12).Value = 9969.231  # L137: 9514.2855 -> 9969.231
$ws.Cells.Item(137, 13).Value = -1076115  # M137: -671991.36 -> -1076115
$ws.Cells.Item(137, 14).Value = -15069.231  # N137: -14614.2855 -> -15069.231

$ws.Cells.Item(138, 8).Value = 2864.1045  # H138: 2880.9849 -> 2864.1045
$ws.Cells.Item(138, 10).Value = 2947.85  # J138: 2968.1526 -> 2947.85
$ws.Cells.Item(138, 12).Value = 8843.549999999999  # L138: 8904.4578 -> 8843.549999999999
$ws.Cells.Item(138, 14).Value = -19123.55  # N138: -19184.4578 -> -19123.55

$ws.Cells.Item(141, 8).Value = 2888  # H141: 2351.25 -> 2888
$ws.Cells.Item(141, 9).Value = 2592.6667  # I141: 2323.5715 -> 2592.6667
$ws.Cells.Item(141, 10).Value = 3183.3333  # J141: 2390 -> 3183.3333
$ws.Cells.Item(141, 11).Value = 7778.000100000001  # K141: 6970.7145 -> 7778.000100000001
$ws.Cells.Item(141, 12).Value = 9549.999899999999  # L141: 7170 -> 9549.999899999999
$ws.Cells.Item(141, 13).Value = -2598.000100000001  # M141: -1790.7145 -> -2598.000100000001
$ws.Cells.Item(141, 14).Value = -19909.9999  # N141: -17530 -> -19909.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3507.7297  # H32: 3516.8784 -> 3507.7297
$ws.Cells.Item(32, 9).Value = 2847.0154  # I32: 2801.9849 -> 2847.0154
$ws.Cells.Item(32, 10).Value = 8279.556  # J32: 9414.75 -> 8279.556
$ws.Cells.Item(32, 11).Value = 2847.0154  # K32: 2801.9849 -> 2847.0154
$ws.Cells.Item(32, 12).Value = 8279.556  # L32: 9414.75 -> 8279.556
$ws.Cells.Item(32, 13).Value = -2560.0154  # M32: -2514.9849 -> -2560.0154
$ws.Cells.Item(32, 14).Value = -8853.556  # N32: -9988.75 -> -8853.556

$ws.Cells.Item(122, 8).Value = 580673.5  # H122: 402241.78 -> 580673.5
$ws.Cells.Item(122, 9).Value = 1668.3793  # I122: 1345.6818 -> 1668.3793
$ws.Cells.Item(122, 10).Value = 2979409  # J122: 2607170.5 -> 2979409
$ws.Cells.Item(122, 11).Value = 5005.1379  # K122: 4037.0454 -> 5005.1379
$ws.Cells.Item(122, 12).Value = 8938227  # L122: 7821511.5 -> 8938227
$ws.Cells.Item(122, 13).Value = -2555.1379  # M122: -1587.0454 -> -2555.1379
$ws.Cells.Item(122, 14).Value = -8943127  # N122: -7826411.5 -> -8943127

$ws.Cells.Item(132, 8).Value = 2689.6  # H132: 2723.4583 -> 2689.6
$ws.Cells.Item(132, 9).Value = 1956.8  # I132: 1958.3334 -> 1956.8
$ws.Cells.Item(132, 10).Value = 3788.8  # J132: 3998.6667 -> 3788.8
$ws.Cells.Item(132, 11).Value = 5870.4  # K132: 5875.0002 -> 5870.4
$ws.Cells.Item(132, 12).Value = 11366.4  # L132: 11996.0001 -> 11366.4
$ws.Cells.Item(132, 13).Value = -3340.4  # M132: -3345.0002 -> -3340.4
$ws.Cells.Item(132, 14).Value = -16426.4  # N132: -17056.0001 -> -16426.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1421.1364  # H20: 1339.0416 -> 1421.1364
$ws.Cells.Item(20, 9).Value = 1174  # I20: 1039.8182 -> 1174
$ws.Cells.Item(20, 11).Value = 1174  # K20: 1039.8182 -> 1174
$ws.Cells.Item(20, 13).Value = -927  # M20: -792.8181999999999 -> -927

$ws.Cells.Item(31, 8).Value = 4000  # H31: 20000 -> 4000
$ws.Cells.Item(31, 9).Value = 5166.6665  # I31: 20000 -> 5166.6665
$ws.Cells.Item(31, 10).Value = 500  # J31: 0 -> 500
$ws.Cells.Item(31, 11).Value = 5166.6665  # K31: 20000 -> 5166.6665
$ws.Cells.Item(31, 12).Value = 500  # L31: 0 -> 500
$ws.Cells.Item(31, 13).Value = -4914.6665  # M31: -19748 -> -4914.6665
$ws.Cells.Item(31, 14).Value = -1004  # N31: None -> -1004

$ws.Cells.Item(94, 8).Value = 10112577  # H94: 11376336 -> 10112577
$ws.Cells.Item(94, 9).Value = 30304530  # I94: 30304696 -> 30304530
$ws.Cells.Item(94, 10).Value = 16600.166  # J94: 19320.2 -> 16600.166
$ws.Cells.Item(94, 11).Value = 30304530  # K94: 30304696 -> 30304530
$ws.Cells.Item(94, 12).Value = 16600.166  # L94: 19320.2 -> 16600.166
$ws.Cells.Item(94, 13).Value = -30304079  # M94: -30304245 -> -30304079
$ws.Cells.Item(94, 14).Value = -17502.166  # N94: -20222.2 -> -17502.166

$ws.Cells.Item(134, 8).Value = 3906.3635  # H134: 3640.5757 -> 3906.3635
$ws.Cells.Item(134, 9).Value = 1769  # I134: 1558.5652 -> 1769
$ws.Cells.Item(134, 10).Value = 6807.0713  # J134: 8429.200000000001 -> 6807.0713
$ws.Cells.Item(134, 11).Value = 5307  # K134: 4675.6956 -> 5307
$ws.Cells.Item(134, 12).Value = 20421.2139  # L134: 25287.6 -> 20421.2139
$ws.Cells.Item(134, 13).Value = -2772  # M134: -2140.6956 -> -2772
$ws.Cells.Item(134, 14).Value = -25491.2139  # N134: -30357.6 -> -25491.2139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 30081.094  # H31: 31753.232 -> 30081.094
$ws.Cells.Item(31, 10).Value = 116349.75  # J31: 153466.67 -> 116349.75
$ws.Cells.Item(31, 12).Value = 116349.75  # L31: 153466.67 -> 116349.75
$ws.Cells.Item(31, 14).Value = -116939.75  # N31: -154056.67 -> -116939.75

$ws.Cells.Item(34, 8).Value = 30081.094  # H34: 31753.232 -> 30081.094
$ws.Cells.Item(34, 10).Value = 116349.75  # J34: 153466.67 -> 116349.75
$ws.Cells.Item(34, 12).Value = 116349.75  # L34: 153466.67 -> 116349.75
$ws.Cells.Item(34, 14).Value = -116753.75  # N34: -153870.67 -> -116753.75

$ws.Cells.Item(58, 8).Value = 1968.9259  # H58: 1934.3334 -> 1968.9259
$ws.Cells.Item(58, 9).Value = 1778.8572  # I58: 1687.6 -> 1778.8572
$ws.Cells.Item(58, 10).Value = 2173.6155  # J58: 2242.75 -> 2173.6155
$ws.Cells.Item(58, 11).Value = 1778.8572  # K58: 1687.6 -> 1778.8572
$ws.Cells.Item(58, 12).Value = 2173.6155  # L58: 2242.75 -> 2173.6155
$ws.Cells.Item(58, 13).Value = -1575.8572  # M58: -1484.6 -> -1575.8572
$ws.Cells.Item(58, 14).Value = -2579.6155  # N58: -2648.75 -> -2579.6155

$ws.Cells.Item(107, 8).Value = 4955.4443  # H107: 4825.125 -> 4955.4443
$ws.Cells.Item(107, 10).Value = 5999  # J107: 5999.5 -> 5999
$ws.Cells.Item(107, 12).Value = 5999  # L107: 5999.5 -> 5999
$ws.Cells.Item(107, 14).Value = -9839  # N107: -9839.5 -> -9839

$ws.Cells.Item(132, 8).Value = 87739.13  # H132: 101110.31 -> 87739.13
$ws.Cells.Item(132, 9).Value = 63980.938  # I132: 72823.21000000001 -> 63980.938
$ws.Cells.Item(132, 10).Value = 151094.33  # J132: 180314.2 -> 151094.33
$ws.Cells.Item(132, 11).Value = 191942.814  # K132: 218469.63 -> 191942.814
$ws.Cells.Item(132, 12).Value = 453282.99  # L132: 540942.6000000001 -> 453282.99
$ws.Cells.Item(132, 13).Value = -189412.814  # M132: -215939.63 -> -189412.814
$ws.Cells.Item(132, 14).Value = -458342.99  # N132: -546002.6000000001 -> -458342.99

$ws.Cells.Item(134, 8).Value = 25696.025  # H134: 25083.879 -> 25696.025
$ws.Cells.Item(134, 9).Value = 39972.207  # I134: 38397.24 -> 39972.207
$ws.Cells.Item(134, 11).Value = 119916.621  # K134: 115191.72 -> 119916.621
$ws.Cells.Item(134, 13).Value = -117381.621  # M134: -112656.72 -> -117381.621

$ws.Cells.Item(136, 8).Value = 1968.9259  # H136: 1934.3334 -> 1968.9259
$ws.Cells.Item(136, 9).Value = 1778.8572  # I136: 1687.6 -> 1778.8572
$ws.Cells.Item(136, 10).Value = 2173.6155  # J136: 2242.75 -> 2173.6155
$ws.Cells.Item(136, 11).Value = 5336.571599999999  # K136: 5062.799999999999 -> 5336.571599999999
$ws.Cells.Item(136, 12).Value = 6520.8465  # L136: 6728.25 -> 6520.8465
$ws.Cells.Item(136, 13).Value = -2786.571599999999  # M136: -2512.799999999999 -> -2786.571599999999
$ws.Cells.Item(136, 14).Value = -11620.8465  # N136: -11828.25 -> -11620.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 604.9231  # H6: 631 -> 604.9231
$ws.Cells.Item(6, 9).Value = 578.5454999999999  # I6: 597.4545000000001 -> 578.5454999999999
$ws.Cells.Item(6, 10).Value = 750  # J6: 1000 -> 750
$ws.Cells.Item(6, 11).Value = 1735.6365  # K6: 1792.3635 -> 1735.6365
$ws.Cells.Item(6, 12).Value = 2250  # L6: 3000 -> 2250
$ws.Cells.Item(6, 13).Value = -1622.6365  # M6: -1679.3635 -> -1622.6365
$ws.Cells.Item(6, 14).Value = -2476  # N6: -3226 -> -2476

$ws.Cells.Item(129, 8).Value = 746.3333  # H129: 2069.7856 -> 746.3333
$ws.Cells.Item(129, 9).Value = 615.6  # I129: 1470.6364 -> 615.6
$ws.Cells.Item(129, 10).Value = 1400  # J129: 4266.6665 -> 1400
$ws.Cells.Item(129, 11).Value = 1846.8  # K129: 4411.9092 -> 1846.8
$ws.Cells.Item(129, 12).Value = 4200  # L129: 12799.9995 -> 4200
$ws.Cells.Item(129, 13).Value = 3153.2  # M129: 588.0907999999999 -> 3153.2
$ws.Cells.Item(129, 14).Value = -14200  # N129: -22799.9995 -> -14200

$ws.Cells.Item(132, 8).Value = 2406.875  # H132: 2431.875 -> 2406.875
$ws.Cells.Item(132, 9).Value = 1725  # I132: 1800 -> 1725
$ws.Cells.Item(132, 11).Value = 15525  # K132: 16200 -> 15525
$ws.Cells.Item(132, 13).Value = -12995  # M132: -13670 -> -12995

$ws.Cells.Item(137, 8).Value = 5979.25  # H137: 7332.3335 -> 5979.25
$ws.Cells.Item(137, 9).Value = 5992.5  # I137: 7999 -> 5992.5
$ws.Cells.Item(137, 10).Value = 5966  # J137: 5999 -> 5966
$ws.Cells.Item(137, 11).Value = 17977.5  # K137: 23997 -> 17977.5
$ws.Cells.Item(137, 12).Value = 17898  # L137: 17997 -> 17898
$ws.Cells.Item(137, 13).Value = -12877.5  # M137: -18897 -> -12877.5
$ws.Cells.Item(137, 14).Value = -28098  # N137: -28197 -> -28098

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 3669.3333  # H29: 3252 -> 3669.3333
$ws.Cells.Item(29, 9).Value = 4500  # I29: 3666.6667 -> 4500
$ws.Cells.Item(29, 11).Value = 4500  # K29: 3666.6667 -> 4500
$ws.Cells.Item(29, 13).Value = -4210  # M29: -3376.6667 -> -4210

$ws.Cells.Item(80, 8).Value = 4883650  # H80: 4070324.8 -> 4883650
$ws.Cells.Item(80, 9).Value = 8131748  # I80: 8132047.5 -> 8131748
$ws.Cells.Item(80, 10).Value = 11503  # J80: 8602 -> 11503
$ws.Cells.Item(80, 11).Value = 8131748  # K80: 8132047.5 -> 8131748
$ws.Cells.Item(80, 12).Value = 11503  # L80: 8602 -> 11503
$ws.Cells.Item(80, 13).Value = -8130750  # M80: -8131049.5 -> -8130750
$ws.Cells.Item(80, 14).Value = -13499  # N80: -10598 -> -13499

$ws.Cells.Item(83, 8).Value = 4883650  # H83: 4070324.8 -> 4883650
$ws.Cells.Item(83, 9).Value = 8131748  # I83: 8132047.5 -> 8131748
$ws.Cells.Item(83, 10).Value = 11503  # J83: 8602 -> 11503
$ws.Cells.Item(83, 11).Value = 40658740  # K83: 40660237.5 -> 40658740
$ws.Cells.Item(83, 12).Value = 57515  # L83: 43010 -> 57515
$ws.Cells.Item(83, 13).Value = -40653748  # M83: -40655245.5 -> -40653748
$ws.Cells.Item(83, 14).Value = -67499  # N83: -52994 -> -67499

$ws.Cells.Item(97, 8).Value = 1254241.2  # H97: 1191543.9 -> 1254241.2
$ws.Cells.Item(97, 9).Value = 1401727.5  # I97: 1323870.1 -> 1401727.5
$ws.Cells.Item(97, 11).Value = 1401727.5  # K97: 1323870.1 -> 1401727.5
$ws.Cells.Item(97, 13).Value = -1401231.5  # M97: -1323374.1 -> -1401231.5

$ws.Cells.Item(122, 8).Value = 389944.25  # H122: 472132.94 -> 389944.25
$ws.Cells.Item(122, 9).Value = 470143.06  # I122: 637787.5600000001 -> 470143.06
$ws.Cells.Item(122, 10).Value = 9000  # J122: 8300 -> 9000
$ws.Cells.Item(122, 11).Value = 1410429.18  # K122: 1913362.68 -> 1410429.18
$ws.Cells.Item(122, 12).Value = 27000  # L122: 24900 -> 27000
$ws.Cells.Item(122, 13).Value = -1407979.18  # M122: -1910912.68 -> -1407979.18
$ws.Cells.Item(122, 14).Value = -31900  # N122: -29800 -> -31900

$ws.Cells.Item(132, 8).Value = 4403  # H132: 3953.5715 -> 4403
$ws.Cells.Item(132, 9).Value = 3174.4443  # I132: 2668.1 -> 3174.4443
$ws.Cells.Item(132, 11).Value = 9523.332900000001  # K132: 8004.299999999999 -> 9523.332900000001
$ws.Cells.Item(132, 13).Value = -6993.332900000001  # M132: -5474.299999999999 -> -6993.332900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3327.125  # H7: 3344.2083 -> 3327.125
$ws.Cells.Item(7, 9).Value = 1865.9375  # I7: 1891.5625 -> 1865.9375
$ws.Cells.Item(7, 11).Value = 1865.9375  # K7: 1891.5625 -> 1865.9375
$ws.Cells.Item(7, 13).Value = -1753.9375  # M7: -1779.5625 -> -1753.9375

$ws.Cells.Item(22, 8).Value = 247722.25  # H22: 165390.67 -> 247722.25
$ws.Cells.Item(22, 9).Value = 247722.25  # I22: 198317.8 -> 247722.25
$ws.Cells.Item(22, 10).Value = 0  # J22: 755 -> 0
$ws.Cells.Item(22, 11).Value = 247722.25  # K22: 198317.8 -> 247722.25
$ws.Cells.Item(22, 12).Value = 0  # L22: 755 -> 0
$ws.Cells.Item(22, 13).ClearContents()  # M22: delete (was -198022.8)
$ws.Cells.Item(22, 14).Value = -247427.25  # N22: -1345 -> -247427.25

$ws.Cells.Item(27, 8).Value = 247722.25  # H27: 165390.67 -> 247722.25
$ws.Cells.Item(27, 9).Value = 247722.25  # I27: 198317.8 -> 247722.25
$ws.Cells.Item(27, 10).Value = 0  # J27: 755 -> 0
$ws.Cells.Item(27, 11).Value = 247722.25  # K27: 198317.8 -> 247722.25
$ws.Cells.Item(27, 12).Value = 0  # L27: 755 -> 0
$ws.Cells.Item(27, 13).ClearContents()  # M27: delete (was -198210.8)
$ws.Cells.Item(27, 14).Value = -247615.25  # N27: -969 -> -247615.25

$ws.Cells.Item(40, 8).Value = 4653.12  # H40: 4805.375 -> 4653.12
$ws.Cells.Item(40, 9).Value = 3385.842  # I40: 3518.4443 -> 3385.842
$ws.Cells.Item(40, 11).Value = 3385.842  # K40: 3518.4443 -> 3385.842
$ws.Cells.Item(40, 13).Value = -3249.842  # M40: -3382.4443 -> -3249.842

$ws.Cells.Item(51, 8).Value = 32979.832  # H51: 42499.75 -> 32979.832
$ws.Cells.Item(51, 10).Value = 32979.832  # J51: 42499.75 -> 32979.832
$ws.Cells.Item(51, 12).Value = 32979.832  # L51: 42499.75 -> 32979.832
$ws.Cells.Item(51, 14).Value = -33935.832  # N51: -43455.75 -> -33935.832

$ws.Cells.Item(61, 8).Value = 3270303.8  # H61: 2926152.5 -> 3270303.8
$ws.Cells.Item(61, 9).Value = 4275910.5  # I61: 4275915 -> 4275910.5
$ws.Cells.Item(61, 10).Value = 2080.875  # J61: 1667.0834 -> 2080.875
$ws.Cells.Item(61, 11).Value = 4275910.5  # K61: 4275915 -> 4275910.5
$ws.Cells.Item(61, 12).Value = 2080.875  # L61: 1667.0834 -> 2080.875
$ws.Cells.Item(61, 13).Value = -4275708.5  # M61: -4275713 -> -4275708.5
$ws.Cells.Item(61, 14).Value = -2484.875  # N61: -2071.0834 -> -2484.875

$ws.Cells.Item(113, 8).Value = 3270303.8  # H113: 2926152.5 -> 3270303.8
$ws.Cells.Item(113, 9).Value = 4275910.5  # I113: 4275915 -> 4275910.5
$ws.Cells.Item(113, 10).Value = 2080.875  # J113: 1667.0834 -> 2080.875
$ws.Cells.Item(113, 11).Value = 4275910.5  # K113: 4275915 -> 4275910.5
$ws.Cells.Item(113, 12).Value = 2080.875  # L113: 1667.0834 -> 2080.875
$ws.Cells.Item(113, 13).Value = -4273740.5  # M113: -4273745 -> -4273740.5
$ws.Cells.Item(113, 14).Value = -6420.875  # N113: -6007.0834 -> -6420.875

$ws.Cells.Item(122, 8).Value = 5440.0454  # H122: 4447.433 -> 5440.0454
$ws.Cells.Item(122, 9).Value = 3645.7273  # I122: 2880.5557 -> 3645.7273
$ws.Cells.Item(122, 10).Value = 7234.364  # J122: 6797.75 -> 7234.364
$ws.Cells.Item(122, 11).Value = 10937.1819  # K122: 8641.667099999999 -> 10937.1819
$ws.Cells.Item(122, 12).Value = 21703.092  # L122: 20393.25 -> 21703.092
$ws.Cells.Item(122, 13).Value = -8487.1819  # M122: -6191.667099999999 -> -8487.1819
$ws.Cells.Item(122, 14).Value = -26603.092  # N122: -25293.25 -> -26603.092

$ws.Cells.Item(126, 8).Value = 3327.125  # H126: 3344.2083 -> 3327.125
$ws.Cells.Item(126, 9).Value = 1865.9375  # I126: 1891.5625 -> 1865.9375
$ws.Cells.Item(126, 11).Value = 5597.8125  # K126: 5674.6875 -> 5597.8125
$ws.Cells.Item(126, 13).Value = -3127.8125  # M126: -3204.6875 -> -3127.8125

$ws.Cells.Item(136, 8).Value = 120659.12  # H136: 170520.25 -> 120659.12
$ws.Cells.Item(136, 9).Value = 128122.81  # I136: 185909.36 -> 128122.81
$ws.Cells.Item(136, 11).Value = 384368.43  # K136: 557728.08 -> 384368.43
$ws.Cells.Item(136, 13).Value = -381818.43  # M136: -555178.08 -> -381818.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2720.4546  # H122: 2008.9445 -> 2720.4546
$ws.Cells.Item(122, 9).Value = 1195  # I122: 1025.25 -> 1195
$ws.Cells.Item(122, 10).Value = 3991.6667  # J122: 2795.9 -> 3991.6667
$ws.Cells.Item(122, 11).Value = 3585  # K122: 3075.75 -> 3585
$ws.Cells.Item(122, 12).Value = 11975.0001  # L122: 8387.700000000001 -> 11975.0001
$ws.Cells.Item(122, 13).Value = -1135  # M122: -625.75 -> -1135
$ws.Cells.Item(122, 14).Value = -16875.0001  # N122: -13287.7 -> -16875.0001

$ws.Cells.Item(136, 8).Value = 2880.158  # H136: 3148.4119 -> 2880.158
$ws.Cells.Item(136, 9).Value = 1659.3572  # I136: 1835.9166 -> 1659.3572
$ws.Cells.Item(136, 11).Value = 4978.071599999999  # K136: 5507.7498 -> 4978.071599999999
$ws.Cells.Item(136, 13).Value = -2428.071599999999  # M136: -2957.7498 -> -2428.071599999999

Write-Output "Applied Hyperion_Profits.xlsx cell updates"
